# Apply the "relational spreadsheet" cutscene_tests.xlsx edit:
#  - Remove the unused END_GAME shared string / row
#  - Shift all cutscene rows up by one (the old leftover "1" in A1 is gone)
#  - Update the "Frog" dialogue row's music cue from "frogs" to "jazzy_retro_battle_theme"
#  - Update the "Frog" dialogue row's choice tag from "0,1" to "LEFT"
#  - Resulting sheet is 5 rows (A1:I5) instead of 7 (A1:I7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gameflow")

# Clear out the old content entirely (rows 1-7, columns A-I) so stale cells don't linger.
$ws.Range("A1:I7").ClearContents()

# Row 1: START_SCENE CUTSCENE jazzy_retro_battle_theme
$ws.Range("A1").Value = "START_SCENE"
$ws.Range("B1").Value = "CUTSCENE"
$ws.Range("C1").Value = "jazzy_retro_battle_theme"

# Row 2: DIALOGUE Tanuki <hey-text> jazzy_retro_battle_theme tanuki_mario 0,1 END_DIALOGUE
$ws.Range("A2").Value = "DIALOGUE"
$ws.Range("B2").Value = "Tanuki"
$ws.Range("C2").Value = "<size=48>Hey you!</size> You're walking in the <color=red>wrong</color> part of town."
$ws.Range("D2").Value = "jazzy_retro_battle_theme"
$ws.Range("E2").Value = "tanuki_mario"
$ws.Range("F2").Value = "0,1"
$ws.Range("G2").Value = "END_DIALOGUE"

# Row 3: DIALOGUE Frog <ribbit-text> jazzy_retro_battle_theme frog_mario LEFT END_DIALOGUE
$ws.Range("A3").Value = "DIALOGUE"
$ws.Range("B3").Value = "Frog"
$ws.Range("C3").Value = "Ribbit <i>Ribbit!</i> <size=24>(Yeah <color=green>frog-face!</color> Wrong part of town!)</size>"
$ws.Range("D3").Value = "jazzy_retro_battle_theme"
$ws.Range("E3").Value = "frog_mario"
$ws.Range("F3").Value = "LEFT"
$ws.Range("G3").Value = "END_DIALOGUE"

# Row 4: DIALOGUE Tanuki "Let's get em!" jazzy_retro_battle_theme tanuki_mario RIGHT frog_mario "1, 1" END_DIALOGUE
$ws.Range("A4").Value = "DIALOGUE"
$ws.Range("B4").Value = "Tanuki"
$ws.Range("C4").Value = "Let's get em!"
$ws.Range("D4").Value = "jazzy_retro_battle_theme"
$ws.Range("E4").Value = "tanuki_mario"
$ws.Range("F4").Value = "RIGHT"
$ws.Range("G4").Value = "frog_mario"
$ws.Range("H4").Value = "1, 1"
$ws.Range("I4").Value = "END_DIALOGUE"

# Row 5: END_SCENE
$ws.Range("A5").Value = "END_SCENE"

# Update the saved selection/active cell like the authored workbook (B11).
$ws.Range("B11").Select()
